$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume data refresh (text-valued cells).
# Force text storage (matches original inline-string cell type) without
# leaving a residual number-format style on the cell.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "316.97"
Set-TextValue "E2" "1.65%"
Set-TextValue "D3" "37.80"
Set-TextValue "E3" "1.14%"
Set-TextValue "D4" "5.159"
Set-TextValue "E4" "0.20%"
Set-TextValue "D5" "0.07970"
Set-TextValue "E5" "1.79%"
Set-TextValue "D6" "4.462"
Set-TextValue "E6" "0.92%"
Set-TextValue "D7" "8.507"
Set-TextValue "E7" "2.84%"
Set-TextValue "D8" "1.930"
Set-TextValue "E8" "1.21%"
Set-TextValue "E9" "6.49%"
Set-TextValue "D10" "0.9398"
Set-TextValue "E10" "2.35%"
Set-TextValue "D11" "0.1273"
Set-TextValue "E11" "7.29%"
Set-TextValue "D12" "0.1937"
Set-TextValue "E12" "0.26%"
Set-TextValue "D13" "0.09021"
Set-TextValue "E13" "-0.86%"
Set-TextValue "D14" "0.03426"
Set-TextValue "E14" "2.10%"
Set-TextValue "D15" "0.09525"
Set-TextValue "E15" "-0.90%"
Set-TextValue "D16" "0.001392"
Set-TextValue "E16" "0.63%"
Set-TextValue "D17" "0.006037"
Set-TextValue "E17" "4.42%"
Set-TextValue "D18" "3.415"
Set-TextValue "E18" "-2.82%"
Set-TextValue "D19" "0.3514"
Set-TextValue "D20" "6.532"
Set-TextValue "E20" "24.21%"
Set-TextValue "D21" "0.1303"
Set-TextValue "E21" "2.41%"
Set-TextValue "D22" "0.2299"
Set-TextValue "E22" "-11.47%"
Set-TextValue "D23" "0.04348"
Set-TextValue "E23" "-0.75%"
Set-TextValue "D24" "0.001197"
Set-TextValue "E24" "-4.34%"
Set-TextValue "D25" "0.004415"
Set-TextValue "E25" "-5.48%"
Set-TextValue "D26" "0.0001323"
Set-TextValue "E26" "-2.94%"
Set-TextValue "D27" "0.0003968"
Set-TextValue "E27" "-0.77%"
Set-TextValue "D39" "0.02368"
Set-TextValue "E39" "3.70%"
Set-TextValue "D40" "0.05174"
Set-TextValue "E40" "2.51%"
Set-TextValue "D41" "0.007422"
Set-TextValue "E41" "-0.18%"
Set-TextValue "D42" "0.1395"
Set-TextValue "E42" "3.59%"
Set-TextValue "D43" "0.008336"
Set-TextValue "E43" "-7.97%"
Set-TextValue "D44" "0.002059"
Set-TextValue "E44" "8.10%"
Set-TextValue "D45" "0.008731"
Set-TextValue "E45" "-7.39%"
Set-TextValue "D46" "0.00006396"
Set-TextValue "E46" "-4.15%"
Set-TextValue "D47" "0.00000000746"
Set-TextValue "E47" "-0.78%"
Set-TextValue "D48" "0.002851"
Set-TextValue "E48" "-13.54%"
Set-TextValue "D49" "0.001680"
Set-TextValue "E49" "67.65%"
Set-TextValue "D50" "0.00002089"
Set-TextValue "E50" "-0.78%"
Set-TextValue "D51" "0.0001989"
Set-TextValue "E51" "-0.78%"
